# Applies Word's automatic proofing markup (<w:proofErr> gramStart/gramEnd
# and spellStart/spellEnd pairs) that appears around certain words once the
# document is re-proofed, splitting the affected runs accordingly.
#
# Because the sandboxed COM layer does not run an actual spell/grammar
# checker, each paragraph/run that needs new <w:proofErr/> markers is
# rebuilt explicitly via Range.InsertXML with the exact WordprocessingML
# fragment Word itself would have produced.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rNs = 'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Get-ParaByText($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "You're implementing a backend DB to the website you created PD1. "
#    -> "You're" gets wrapped in gramStart/gramEnd and split into its own run.
# ---------------------------------------------------------------------
$p1 = Get-ParaByText $d "implementing a backend DB"
$r1 = $p1.Range
$apo = [char]0x2019
$xml1 = "<w:p $wNs $w14Ns w14:paraId=`"2287CDB0`" w14:textId=`"77777777`" w:rsidR=`"00070F5B`" w:rsidRDefault=`"00070F5B`" w:rsidP=`"00070F5B`"><w:proofErr w:type=`"gramStart`"/><w:r><w:t>You${apo}re</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> implementing a backend DB to the website you created PD1. </w:t></w:r></w:p>"
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "Develop a Relational diagram ... are SmartDraw, draw.io, erdplus, etc."
#    -> "SmartDraw" and "erdplus" each get wrapped in spellStart/spellEnd.
# ---------------------------------------------------------------------
$p2 = Get-ParaByText $d "Develop a Relational diagram"
$r2 = $p2.Range
$xml2 = @"
<w:p $wNs $w14Ns w14:paraId="1A3E6BF9" w14:textId="77777777" w:rsidR="00070F5B" w:rsidRDefault="00070F5B" w:rsidP="00070F5B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Develop a Relational diagram with all the constraints and provide a table showing the data types, such as string, number, etc. Free diagramming tools (feel free to use and/or share any other that you find) are </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SmartDraw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, draw.io, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>erdplus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, etc.</w:t></w:r></w:p>
"@
$r2.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) "IMPACT is a non-profit organization with the officers, volunteers and events."
#    -> "volunteers" gets wrapped in gramStart/gramEnd.
# ---------------------------------------------------------------------
$p3 = Get-ParaByText $d "IMPACT is a non-profit organization"
$r3 = $p3.Range
$xml3 = "<w:p $wNs $w14Ns w14:paraId=`"3582F904`" w14:textId=`"77777777`" w:rsidR=`"00070F5B`" w:rsidRDefault=`"00070F5B`" w:rsidP=`"00070F5B`"><w:r><w:t xml:space=`"preserve`">IMPACT is a non-profit organization with the officers, </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>volunteers</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> and events.</w:t></w:r></w:p>"
$r3.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) "There are many free diagramming tools available such as SmartDraw
#    (www.smartdraw.com), and draw.io" -> the already-italicised "SmartDraw"
#    run gets wrapped in spellStart/spellEnd, formatting/hyperlink untouched.
# ---------------------------------------------------------------------
$p4 = Get-ParaByText $d "There are many free diagramming tools available"
$full4 = $p4.Range
$find4 = $full4.Duplicate
$find4.Find.Execute("SmartDraw", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4 = $d.Range($find4.Start, $full4.End)
$xml4 = @"
<w:p $wNs $rNs><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009262B0"><w:rPr><w:i/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>SmartDraw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009262B0"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:hyperlink r:id="rId6" w:history="1"><w:r w:rsidRPr="009262B0"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>www.smartdraw.com</w:t></w:r></w:hyperlink><w:r w:rsidRPr="009262B0"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">), and </w:t></w:r><w:r w:rsidRPr="009262B0"><w:rPr><w:i/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>draw.io</w:t></w:r></w:p>
"@
$r4.InsertXML($xml4)

Write-Output "done"
